$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.640.01'
$ws.Range("E2").Value = '  +3.13%  '
$ws.Range("D3").Value = '1.852.79'
$ws.Range("D4").Value = '''1.030'
$ws.Range("E4").Value = '  +2.44%  '
$ws.Range("D5").Value = '''321.62'
$ws.Range("E5").Value = '  +4.06%  '
$ws.Range("D6").Value = '''1.029'
$ws.Range("E6").Value = '  +2.43%  '
$ws.Range("D7").Value = '''0.4380'
$ws.Range("E7").Value = '  +1.40%  '
$ws.Range("D8").Value = '''0.3757'
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("D9").Value = '''0.07419'
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = '''0.8765'
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("D11").Value = '''21.50'
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("D12").Value = '1.869.26'
$ws.Range("E12").Value = '  -6.13%  '
$ws.Range("D14").Value = '''6.705'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '''0.07190'
$ws.Range("E15").Value = '  +4.00%  '
$ws.Range("D16").Value = '''82.98'
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("D17").Value = '''1.035'
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("D18").Value = '''0.000009060'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '''1.029'
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").Value = '''15.47'
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("D21").Value = '27.654.71'
$ws.Range("E21").Value = '  +3.00%  '
$ws.Range("D22").Value = '''5.274'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").Value = '''11.25'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '2.078.07'
$ws.Range("E24").Value = '  -6.39%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''1.948'
$ws.Range("E25").Value = '  +3.96%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''157.53'
$ws.Range("E26").Value = '  +2.53%  '
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("D28").Value = '''5.300'
$ws.Range("E28").Value = '  +1.93%  '
$ws.Range("D29").Value = '''1.937'
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("D30").Value = '''116.40'
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("E32").Value = '  +3.79%  '
$ws.Range("D33").Value = '''0.7700'
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D34").Value = '''4.523'
$ws.Range("E34").Value = '  +1.98%  '
$ws.Range("D35").Value = '''2.882'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").Value = '''1.031'
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("D37").Value = '''1.154'
$ws.Range("E37").Value = '  +2.51%  '
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("D39").Value = '''0.05290'
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("D40").Value = '''2.832'
$ws.Range("E40").Value = '  +6.25%  '
$ws.Range("D41").Value = '''0.5189'
$ws.Range("E41").Value = '  +2.23%  '
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").Value = '''6.742'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").Value = '''8.591'
$ws.Range("E44").Value = '  +3.59%  '
$ws.Range("D45").Value = '''108.97'
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("E46").Value = '  +2.39%  '
$ws.Range("E47").Value = '  +4.07%  '
$ws.Range("D48").Value = '''0.4665'
$ws.Range("E48").Value = '  +2.42%  '
$ws.Range("D49").Value = '''0.06392'
$ws.Range("E49").Value = '  +1.79%  '
$ws.Range("D50").Value = '''1.891'
$ws.Range("E50").Value = '  +5.00%  '
$ws.Range("D51").Value = '''39.61'
$ws.Range("E51").Value = '  +5.74%  '
